$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new data rows (real dataset windows) below the existing data
$ws.Range("A22").Value = 998
$ws.Range("B22").Value = 5
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 0.53865740740740742
$ws.Range("D22").NumberFormat = "h:mm:ss"

$ws.Range("A23").Value = 991
$ws.Range("B23").Value = 6
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 0.53865740740740742
$ws.Range("D23").NumberFormat = "h:mm:ss"

# Update the active selection to reflect where the user would continue entering data
$ws.Range("A24").Select()
